$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.861.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.77%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.860.27"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.32%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6357"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.47%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2997"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07475"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.47"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07685"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.866.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.050"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6887"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "84.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.55%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009379"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.093"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "29.845.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.116.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.71%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "238.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.15%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.355"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.30"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1418"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.575"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.99%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.505"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.57%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06087"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.268"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.140"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.143"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.866"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.18%  "
$ws.Range("E35").Value = "  +3.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7280"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.620"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.20%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.856"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01794"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.220.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9292"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.317"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.024.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "66.35"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000124"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5095"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.297"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4087"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.1142"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.22%  "
